$d = $word.ActiveDocument

# 1) Insert four new empty paragraphs (matching the formatting of the
#    existing trailing empty paragraph) right before the end of the
#    document / before sectPr. Assigning a single carriage return to a
#    zero-length Range positioned right before the very end of the
#    story produces a clean empty paragraph with no phantom run, and
#    inherits the paragraph mark formatting (spacing/jc/rPr) of the
#    paragraph that used to be last.
for ($i = 0; $i -lt 4; $i++) {
    $endPos = $d.Content.End
    $insPoint = $d.Range($endPos - 1, $endPos - 1)
    $insPoint.Text = "`r"
}

# 2) Move the "_GoBack" bookmark from its original location (end of the
#    " program matrikulasi." paragraph) onto the new last paragraph.
#    Re-adding a bookmark named "_GoBack" relocates it (removing the
#    previous one). The target Range must span across the paragraph
#    boundary (include the end of the prior paragraph through the end
#    of the new last paragraph) to land correctly on the last
#    paragraph.
$lastPara = $d.Paragraphs.Last
$bmRange = $d.Range($lastPara.Range.Start - 1, $lastPara.Range.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
